# Scheduled market-data refresh: updates Universalis price snapshots
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the per-job leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 176  # H33 currentAveragePrice
$ws.Cells.Item(33, 9).Value = 179.6  # I33 currentAveragePriceNQ
$ws.Cells.Item(33, 10).Value = 149  # J33 currentAveragePriceHQ
$ws.Cells.Item(33, 11).Value = 179.6  # K33 LevePriceNQ
$ws.Cells.Item(33, 12).Value = 149  # L33 LevePriceHQ
$ws.Cells.Item(33, 13).Value = 49.40000000000001  # M33 LeveProfitNQ
$ws.Cells.Item(33, 14).Value = -607  # N33 LeveProfitHQ

$ws.Cells.Item(40, 8).Value = 4369.7  # H40 currentAveragePrice
$ws.Cells.Item(40, 10).Value = 3672.5715  # J40 currentAveragePriceHQ
$ws.Cells.Item(40, 12).Value = 3672.5715  # L40 LevePriceHQ
$ws.Cells.Item(40, 14).Value = -4022.5715  # N40 LeveProfitHQ

$ws.Cells.Item(43, 8).Value = 6419.077  # H43 currentAveragePrice
$ws.Cells.Item(43, 10).Value = 2250  # J43 currentAveragePriceHQ
$ws.Cells.Item(43, 12).Value = 2250  # L43 LevePriceHQ
$ws.Cells.Item(43, 14).Value = -2388  # N43 LeveProfitHQ

$ws.Cells.Item(86, 8).Value = 62503236  # H86 currentAveragePrice
$ws.Cells.Item(86, 9).Value = 90911520  # I86 currentAveragePriceNQ
$ws.Cells.Item(86, 11).Value = 90911520  # K86 LevePriceNQ
$ws.Cells.Item(86, 13).Value = -90910397  # M86 LeveProfitNQ

$ws.Cells.Item(89, 8).Value = 62503236  # H89 currentAveragePrice
$ws.Cells.Item(89, 9).Value = 90911520  # I89 currentAveragePriceNQ
$ws.Cells.Item(89, 11).Value = 454557600  # K89 LevePriceNQ
$ws.Cells.Item(89, 13).Value = -454551984  # M89 LeveProfitNQ

$ws.Cells.Item(132, 8).Value = 9323.689  # H132 currentAveragePrice
$ws.Cells.Item(132, 9).Value = 6625.826  # I132 currentAveragePriceNQ
$ws.Cells.Item(132, 10).Value = 19665.5  # J132 currentAveragePriceHQ
$ws.Cells.Item(132, 11).Value = 19877.478  # K132 LevePriceNQ
$ws.Cells.Item(132, 12).Value = 58996.5  # L132 LevePriceHQ
$ws.Cells.Item(132, 13).Value = -17347.478  # M132 LeveProfitNQ
$ws.Cells.Item(132, 14).Value = -64056.5  # N132 LeveProfitHQ

$ws.Cells.Item(137, 8).Value = 2278316.8  # H137 currentAveragePrice
$ws.Cells.Item(137, 9).Value = 2500898.2  # I137 currentAveragePriceNQ
$ws.Cells.Item(137, 11).Value = 7502694.600000001  # K137 LevePriceNQ
$ws.Cells.Item(137, 13).Value = -7500144.600000001  # M137 LeveProfitNQ

$ws.Cells.Item(141, 8).Value = 5812.1  # H141 currentAveragePrice
$ws.Cells.Item(141, 10).Value = 12600  # J141 currentAveragePriceHQ
$ws.Cells.Item(141, 12).Value = 37800  # L141 LevePriceHQ
$ws.Cells.Item(141, 14).Value = -48160  # N141 LeveProfitHQ

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 7449  # H13 currentAveragePrice
$ws.Cells.Item(13, 9).Value = 7449  # I13 currentAveragePriceNQ
$ws.Cells.Item(13, 11).Value = 7449  # K13 LevePriceNQ
$ws.Cells.Item(13, 13).Value = -7305  # M13 LeveProfitNQ

$ws.Cells.Item(32, 8).Value = 4584.7676  # H32 currentAveragePrice
$ws.Cells.Item(32, 10).Value = 7617  # J32 currentAveragePriceHQ
$ws.Cells.Item(32, 12).Value = 7617  # L32 LevePriceHQ
$ws.Cells.Item(32, 14).Value = -8191  # N32 LeveProfitHQ

$ws.Cells.Item(45, 8).Value = 27209.117  # H45 currentAveragePrice
$ws.Cells.Item(45, 9).Value = 31752  # I45 currentAveragePriceNQ
$ws.Cells.Item(45, 10).Value = 6009  # J45 currentAveragePriceHQ
$ws.Cells.Item(45, 11).Value = 31752  # K45 LevePriceNQ
$ws.Cells.Item(45, 12).Value = 6009  # L45 LevePriceHQ
$ws.Cells.Item(45, 13).Value = -31375  # M45 LeveProfitNQ
$ws.Cells.Item(45, 14).Value = -6763  # N45 LeveProfitHQ

$ws.Cells.Item(74, 8).Value = 329054.47  # H74 currentAveragePrice
$ws.Cells.Item(74, 9).Value = 464549.66  # I74 currentAveragePriceNQ
$ws.Cells.Item(74, 11).Value = 464549.66  # K74 LevePriceNQ
$ws.Cells.Item(74, 13).Value = -463675.66  # M74 LeveProfitNQ

$ws.Cells.Item(77, 8).Value = 329054.47  # H77 currentAveragePrice
$ws.Cells.Item(77, 9).Value = 464549.66  # I77 currentAveragePriceNQ
$ws.Cells.Item(77, 11).Value = 2322748.3  # K77 LevePriceNQ
$ws.Cells.Item(77, 13).Value = -2318380.3  # M77 LeveProfitNQ

$ws.Cells.Item(110, 8).Value = 2938.652  # H110 currentAveragePrice
$ws.Cells.Item(110, 9).Value = 1600.875  # I110 currentAveragePriceNQ
$ws.Cells.Item(110, 10).Value = 5996.4287  # J110 currentAveragePriceHQ
$ws.Cells.Item(110, 11).Value = 1600.875  # K110 LevePriceNQ
$ws.Cells.Item(110, 12).Value = 5996.4287  # L110 LevePriceHQ
$ws.Cells.Item(110, 13).Value = 444.125  # M110 LeveProfitNQ
$ws.Cells.Item(110, 14).Value = -10086.4287  # N110 LeveProfitHQ

$ws.Cells.Item(122, 8).Value = 2999.8572  # H122 currentAveragePrice
$ws.Cells.Item(122, 10).Value = 4825.2  # J122 currentAveragePriceHQ
$ws.Cells.Item(122, 12).Value = 14475.6  # L122 LevePriceHQ
$ws.Cells.Item(122, 14).Value = -19375.6  # N122 LeveProfitHQ

$ws.Cells.Item(132, 8).Value = 1820.8723  # H132 currentAveragePrice
$ws.Cells.Item(132, 9).Value = 1003.6061  # I132 currentAveragePriceNQ
$ws.Cells.Item(132, 10).Value = 3747.2856  # J132 currentAveragePriceHQ
$ws.Cells.Item(132, 11).Value = 3010.8183  # K132 LevePriceNQ
$ws.Cells.Item(132, 12).Value = 11241.8568  # L132 LevePriceHQ
$ws.Cells.Item(132, 13).Value = -480.8182999999999  # M132 LeveProfitNQ
$ws.Cells.Item(132, 14).Value = -16301.8568  # N132 LeveProfitHQ

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3737.41  # H134 currentAveragePrice
$ws.Cells.Item(134, 9).Value = 3570.06  # I134 currentAveragePriceNQ
$ws.Cells.Item(134, 10).Value = 4498.091  # J134 currentAveragePriceHQ
$ws.Cells.Item(134, 11).Value = 10710.18  # K134 LevePriceNQ
$ws.Cells.Item(134, 12).Value = 13494.273  # L134 LevePriceHQ
$ws.Cells.Item(134, 13).Value = -8175.18  # M134 LeveProfitNQ
$ws.Cells.Item(134, 14).Value = -18564.273  # N134 LeveProfitHQ

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 70430.4  # H59 currentAveragePrice
$ws.Cells.Item(59, 10).Value = 79050.664  # J59 currentAveragePriceHQ
$ws.Cells.Item(59, 12).Value = 79050.664  # L59 LevePriceHQ
$ws.Cells.Item(59, 14).Value = -81340.664  # N59 LeveProfitHQ

$ws.Cells.Item(122, 8).Value = 3157.2856  # H122 currentAveragePrice
$ws.Cells.Item(122, 9).Value = 2604.4666  # I122 currentAveragePriceNQ
$ws.Cells.Item(122, 11).Value = 7813.399800000001  # K122 LevePriceNQ
$ws.Cells.Item(122, 13).Value = -5363.399800000001  # M122 LeveProfitNQ

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 8480  # H62 currentAveragePrice
$ws.Cells.Item(62, 10).Value = 8480  # J62 currentAveragePriceHQ
$ws.Cells.Item(62, 12).Value = 25440  # L62 LevePriceHQ
$ws.Cells.Item(62, 14).Value = -26812  # N62 LeveProfitHQ

$ws.Cells.Item(65, 8).Value = 8480  # H65 currentAveragePrice
$ws.Cells.Item(65, 10).Value = 8480  # J65 currentAveragePriceHQ
$ws.Cells.Item(65, 12).Value = 76320  # L65 LevePriceHQ
$ws.Cells.Item(65, 14).Value = -83184  # N65 LeveProfitHQ

$ws.Cells.Item(97, 8).Value = 1300373  # H97 currentAveragePrice
$ws.Cells.Item(97, 9).Value = 5000000  # I97 currentAveragePriceNQ
$ws.Cells.Item(97, 10).Value = 67164  # J97 currentAveragePriceHQ
$ws.Cells.Item(97, 11).Value = 15000000  # K97 LevePriceNQ
$ws.Cells.Item(97, 12).Value = 201492  # L97 LevePriceHQ
$ws.Cells.Item(97, 13).Value = -14999504  # M97 LeveProfitNQ
$ws.Cells.Item(97, 14).Value = -202484  # N97 LeveProfitHQ

$ws.Cells.Item(116, 8).Value = 1212972.4  # H116 currentAveragePrice
$ws.Cells.Item(116, 9).Value = 1816992.1  # I116 currentAveragePriceNQ
$ws.Cells.Item(116, 10).Value = 4932.6665  # J116 currentAveragePriceHQ
$ws.Cells.Item(116, 11).Value = 5450976.300000001  # K116 LevePriceNQ
$ws.Cells.Item(116, 12).Value = 14797.9995  # L116 LevePriceHQ
$ws.Cells.Item(116, 13).Value = -5447534.300000001  # M116 LeveProfitNQ
$ws.Cells.Item(116, 14).Value = -21681.9995  # N116 LeveProfitHQ

$ws.Cells.Item(131, 8).Value = 2440.3704  # H131 currentAveragePrice
$ws.Cells.Item(131, 9).Value = 2484.5  # I131 currentAveragePriceNQ
$ws.Cells.Item(131, 10).Value = 2405.0667  # J131 currentAveragePriceHQ
$ws.Cells.Item(131, 11).Value = 7453.5  # K131 LevePriceNQ
$ws.Cells.Item(131, 12).Value = 7215.2001  # L131 LevePriceHQ
$ws.Cells.Item(131, 13).Value = -2413.5  # M131 LeveProfitNQ
$ws.Cells.Item(131, 14).Value = -17295.2001  # N131 LeveProfitHQ

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2999.8  # H102 currentAveragePrice
$ws.Cells.Item(102, 9).Value = 2999.6667  # I102 currentAveragePriceNQ
$ws.Cells.Item(102, 11).Value = 2999.6667  # K102 LevePriceNQ
$ws.Cells.Item(102, 13).Value = -1377.6667  # M102 LeveProfitNQ

$ws.Cells.Item(132, 8).Value = 3489.3333  # H132 currentAveragePrice
$ws.Cells.Item(132, 10).Value = 3638.0908  # J132 currentAveragePriceHQ
$ws.Cells.Item(132, 12).Value = 10914.2724  # L132 LevePriceHQ
$ws.Cells.Item(132, 14).Value = -15974.2724  # N132 LeveProfitHQ

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1050  # H22 currentAveragePrice
$ws.Cells.Item(22, 9).Value = 750  # I22 currentAveragePriceNQ
$ws.Cells.Item(22, 10).Value = 1150  # J22 currentAveragePriceHQ
$ws.Cells.Item(22, 11).Value = 750  # K22 LevePriceNQ
$ws.Cells.Item(22, 12).Value = 1150  # L22 LevePriceHQ
$ws.Cells.Item(22, 13).Value = -455  # M22 LeveProfitNQ
$ws.Cells.Item(22, 14).Value = -1740  # N22 LeveProfitHQ

$ws.Cells.Item(27, 8).Value = 1050  # H27 currentAveragePrice
$ws.Cells.Item(27, 9).Value = 750  # I27 currentAveragePriceNQ
$ws.Cells.Item(27, 10).Value = 1150  # J27 currentAveragePriceHQ
$ws.Cells.Item(27, 11).Value = 750  # K27 LevePriceNQ
$ws.Cells.Item(27, 12).Value = 1150  # L27 LevePriceHQ
$ws.Cells.Item(27, 13).Value = -643  # M27 LeveProfitNQ
$ws.Cells.Item(27, 14).Value = -1364  # N27 LeveProfitHQ

$ws.Cells.Item(46, 8).Value = 1519.8  # H46 currentAveragePrice
$ws.Cells.Item(46, 10).Value = 1519.8  # J46 currentAveragePriceHQ
$ws.Cells.Item(46, 12).Value = 1519.8  # L46 LevePriceHQ
$ws.Cells.Item(46, 14).Value = -1895.8  # N46 LeveProfitHQ

$ws.Cells.Item(100, 8).Value = 2291.8333  # H100 currentAveragePrice
$ws.Cells.Item(100, 9).Value = 2168  # I100 currentAveragePriceNQ
$ws.Cells.Item(100, 10).Value = 2539.5  # J100 currentAveragePriceHQ
$ws.Cells.Item(100, 11).Value = 2168  # K100 LevePriceNQ
$ws.Cells.Item(100, 12).Value = 2539.5  # L100 LevePriceHQ
$ws.Cells.Item(100, 13).Value = -1627  # M100 LeveProfitNQ
$ws.Cells.Item(100, 14).Value = -3621.5  # N100 LeveProfitHQ

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(82, 8).Value = 60301  # H82 currentAveragePrice
$ws.Cells.Item(82, 10).Value = 60301  # J82 currentAveragePriceHQ
$ws.Cells.Item(82, 12).Value = 60301  # L82 LevePriceHQ
$ws.Cells.Item(82, 14).Value = -61067  # N82 LeveProfitHQ

$ws.Cells.Item(85, 8).Value = 60301  # H85 currentAveragePrice
$ws.Cells.Item(85, 10).Value = 60301  # J85 currentAveragePriceHQ
$ws.Cells.Item(85, 12).Value = 60301  # L85 LevePriceHQ
$ws.Cells.Item(85, 14).Value = -62953  # N85 LeveProfitHQ
